$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the surrounding double-quote characters from the name (column C)
# values for rows 1-10.
for ($r = 1; $r -le 10; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $cCell.Value().Replace('"', '')
}

# Strip the surrounding double-quote characters from the email (column E)
# values for rows 1-10.
for ($r = 1; $r -le 10; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $eCell.Value().Replace('"', '')
}

# Update the active selection to E9, matching the saved view state.
$ws.Range("E9").Select()
